# Apply the two kinds of changes described by the commit:
#  1. The "datetimeFigureOut" date field cached on the slide master and on
#     every slide layout is bumped from 12/29/2020 to 2/8/2021.
#  2. The three mitigation-stage rectangle labels on the slide are renamed
#     from the old "*processing" terminology to the new "*-estimator
#     mitigation" terminology.

$p = $ppt.ActivePresentation

$oldDate = "12/29/2020"
$newDate = "2/8/2021"

# --- 1. Update the date placeholder text everywhere it is cached --------

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DateShapes $master.Shapes

# Every slide layout's date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}

# --- 2. Rename the mitigation-stage rectangles ---------------------------

$renames = @{
    "preprocessing"  = "pre-estimator mitigation"
    "inprocessing"   = "in-estimator mitigation"
    "postprocessing" = "post-estimator mitigation"
}

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            $cur = $tr.Text
            if ($renames.ContainsKey($cur)) {
                $tr.Text = $renames[$cur]
            }
        }
    }
}
